$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.002991666666666667
$ws.Range("H2").Value = 0.008975
$ws.Range("I2").Value = 0.0003566413595017623
$ws.Range("J2").Value = 0.0003566413595017623
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5503976666666667
$ws.Range("N2").Value = 1.651193
$ws.Range("O2").Value = 0.8380541753160734
$ws.Range("P2").Value = 0.8380541753160734
$ws.Range("Q2").Value = 0.001646606352777778
$ws.Range("R2").Value = 0.014819457175
$ws.Range("S2").Value = 0.0002988847804208527
$ws.Range("T2").Value = 0.0002988847804208527
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.002991666666666667
$ws.Range("H3").Value = 0.008975
$ws.Range("I3").Value = 0.0003566413595017623
$ws.Range("J3").Value = 0.0003566413595017623
$ws.Range("O3").Value = 0.1511828328097164
$ws.Range("P3").Value = 0.1511828328097164
$ws.Range("Q3").Value = 0.0002970435805555555
$ws.Range("R3").Value = 0.002673392225
$ws.Range("S3").Value = 0.0000539180510265849
$ws.Range("T3").Value = 0.0000539180510265849
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.002991666666666667
$ws.Range("H4").Value = 0.008975
$ws.Range("I4").Value = 0.0003566413595017623
$ws.Range("J4").Value = 0.0003566413595017623
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.007068666666666667
$ws.Range("N4").Value = 0.021206
$ws.Range("O4").Value = 0.01076299187421013
$ws.Range("P4").Value = 0.01076299187421013
$ws.Range("Q4").Value = 0.00002114709444444445
$ws.Range("R4").Value = 0.00019032385
$ws.Range("S4").Value = 0.000003838528054324723
$ws.Range("T4").Value = 0.000003838528054324723
$ws.Range("I5").Value = 0.9971069332391614
$ws.Range("J5").Value = 0.9971069332391616
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.5503976666666667
$ws.Range("N5").Value = 1.651193
$ws.Range("O5").Value = 0.8380541753160734
$ws.Range("P5").Value = 0.8380541753160734
$ws.Range("Q5").Value = 4.603623687852888
$ws.Range("R5").Value = 41.432613190676
$ws.Range("S5").Value = 0.8356296286376845
$ws.Range("T5").Value = 0.8356296286376846
$ws.Range("I6").Value = 0.9971069332391614
$ws.Range("J6").Value = 0.9971069332391616
$ws.Range("O6").Value = 0.1511828328097164
$ws.Range("P6").Value = 0.1511828328097164
$ws.Range("S6").Value = 0.1507454507813052
$ws.Range("T6").Value = 0.1507454507813052
$ws.Range("I7").Value = 0.9971069332391614
$ws.Range("J7").Value = 0.9971069332391616
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.007068666666666667
$ws.Range("N7").Value = 0.021206
$ws.Range("O7").Value = 0.01076299187421013
$ws.Range("P7").Value = 0.01076299187421013
$ws.Range("Q7").Value = 0.05912358151022221
$ws.Range("R7").Value = 0.532112233592
$ws.Range("S7").Value = 0.01073185382017168
$ws.Range("T7").Value = 0.01073185382017168
$ws.Range("G8").Value = 0.02127666666666667
$ws.Range("H8").Value = 0.06383
$ws.Range("I8").Value = 0.002536425401336767
$ws.Range("J8").Value = 0.002536425401336767
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5503976666666667
$ws.Range("N8").Value = 1.651193
$ws.Range("O8").Value = 0.8380541753160734
$ws.Range("P8").Value = 0.8380541753160734
$ws.Range("Q8").Value = 0.01171062768777778
$ws.Range("R8").Value = 0.10539564919
$ws.Range("S8").Value = 0.002125661897968025
$ws.Range("T8").Value = 0.002125661897968025
$ws.Range("G9").Value = 0.02127666666666667
$ws.Range("H9").Value = 0.06383
$ws.Range("I9").Value = 0.002536425401336767
$ws.Range("J9").Value = 0.002536425401336767
$ws.Range("O9").Value = 0.1511828328097164
$ws.Range("P9").Value = 0.1511828328097164
$ws.Range("Q9").Value = 0.002112567325555555
$ws.Range("R9").Value = 0.01901310593
$ws.Range("S9").Value = 0.0003834639773846144
$ws.Range("T9").Value = 0.0003834639773846144
$ws.Range("G10").Value = 0.02127666666666667
$ws.Range("H10").Value = 0.06383
$ws.Range("I10").Value = 0.002536425401336767
$ws.Range("J10").Value = 0.002536425401336767
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.007068666666666667
$ws.Range("N10").Value = 0.021206
$ws.Range("O10").Value = 0.01076299187421013
$ws.Range("P10").Value = 0.01076299187421013
$ws.Range("Q10").Value = 0.00002114709444444445
$ws.Range("R10").Value = 0.00019032385
$ws.Range("S10").Value = 0.0000272995259841278
$ws.Range("T10").Value = 0.0000272995259841278
